$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Content.Find.Execute("2025-10-11 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-12 Sunday", 2) | Out-Null

# Update each arithmetic expression cell in the practice table
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "87-27="
$t.Cell(1,2).Range.Text = "60-60="
$t.Cell(1,3).Range.Text = "39+1="
$t.Cell(1,4).Range.Text = "92-51="
$t.Cell(1,5).Range.Text = "97-65="
$t.Cell(2,1).Range.Text = "28+6="
$t.Cell(2,2).Range.Text = "84-70="
$t.Cell(2,3).Range.Text = "80-3="
$t.Cell(2,4).Range.Text = "90-8="
$t.Cell(2,5).Range.Text = "13+81="
$t.Cell(3,1).Range.Text = "75-64="
$t.Cell(3,2).Range.Text = "47+0="
$t.Cell(3,3).Range.Text = "61-22="
$t.Cell(3,4).Range.Text = "64+10="
$t.Cell(3,5).Range.Text = "39+10="
$t.Cell(4,1).Range.Text = "85-25="
$t.Cell(4,2).Range.Text = "1+13="
$t.Cell(4,3).Range.Text = "12-2="
$t.Cell(4,4).Range.Text = "58+18="
$t.Cell(4,5).Range.Text = "70+28="
$t.Cell(5,1).Range.Text = "81-75="
$t.Cell(5,2).Range.Text = "55-50="
$t.Cell(5,3).Range.Text = "17+20="
$t.Cell(5,4).Range.Text = "68+9="
$t.Cell(5,5).Range.Text = "28+63="
$t.Cell(6,1).Range.Text = "82-79="
$t.Cell(6,2).Range.Text = "93-14="
$t.Cell(6,3).Range.Text = "16+52="
$t.Cell(6,4).Range.Text = "79-52="
$t.Cell(6,5).Range.Text = "47+50="
$t.Cell(7,1).Range.Text = "26-20="
$t.Cell(7,2).Range.Text = "52+41="
$t.Cell(7,3).Range.Text = "60+25="
$t.Cell(7,4).Range.Text = "31-10="
$t.Cell(7,5).Range.Text = "14+48="
$t.Cell(8,1).Range.Text = "2+34="
$t.Cell(8,2).Range.Text = "88-1="
$t.Cell(8,3).Range.Text = "15-15="
$t.Cell(8,4).Range.Text = "62+27="
$t.Cell(8,5).Range.Text = "17+16="
$t.Cell(9,1).Range.Text = "6+32="
$t.Cell(9,2).Range.Text = "79-65="
$t.Cell(9,3).Range.Text = "52+18="
$t.Cell(9,4).Range.Text = "36+56="
$t.Cell(9,5).Range.Text = "17+68="
$t.Cell(10,1).Range.Text = "15-14="
$t.Cell(10,2).Range.Text = "62-11="
$t.Cell(10,3).Range.Text = "51-12="
$t.Cell(10,4).Range.Text = "4+61="
$t.Cell(10,5).Range.Text = "83-41="
$t.Cell(11,1).Range.Text = "4+33="
$t.Cell(11,2).Range.Text = "27-23="
$t.Cell(11,3).Range.Text = "51-17="
$t.Cell(11,4).Range.Text = "18+50="
$t.Cell(11,5).Range.Text = "87+9="
$t.Cell(12,1).Range.Text = "99-76="
$t.Cell(12,2).Range.Text = "92-86="
$t.Cell(12,3).Range.Text = "29+19="
$t.Cell(12,4).Range.Text = "89+2="
$t.Cell(12,5).Range.Text = "99-25="
$t.Cell(13,1).Range.Text = "99-21="
$t.Cell(13,2).Range.Text = "47+36="
$t.Cell(13,3).Range.Text = "2+39="
$t.Cell(13,4).Range.Text = "24+51="
$t.Cell(13,5).Range.Text = "22-8="
$t.Cell(14,1).Range.Text = "36+1="
$t.Cell(14,2).Range.Text = "36-30="
$t.Cell(14,3).Range.Text = "25-18="
$t.Cell(14,4).Range.Text = "25+21="
$t.Cell(14,5).Range.Text = "90+5="
$t.Cell(15,1).Range.Text = "61+35="
$t.Cell(15,2).Range.Text = "22+48="
$t.Cell(15,3).Range.Text = "53+30="
$t.Cell(15,4).Range.Text = "49+4="
$t.Cell(15,5).Range.Text = "86-56="
$t.Cell(16,1).Range.Text = "25+58="
$t.Cell(16,2).Range.Text = "50+33="
$t.Cell(16,3).Range.Text = "11+87="
$t.Cell(16,4).Range.Text = "56+38="
$t.Cell(16,5).Range.Text = "48-34="
$t.Cell(17,1).Range.Text = "27+2="
$t.Cell(17,2).Range.Text = "15+0="
$t.Cell(17,3).Range.Text = "22+35="
$t.Cell(17,4).Range.Text = "88-60="
$t.Cell(17,5).Range.Text = "41+38="
$t.Cell(18,1).Range.Text = "15+37="
$t.Cell(18,2).Range.Text = "68-32="
$t.Cell(18,3).Range.Text = "58+14="
$t.Cell(18,4).Range.Text = "8+46="
$t.Cell(18,5).Range.Text = "30+3="
$t.Cell(19,1).Range.Text = "18+11="
$t.Cell(19,2).Range.Text = "76-63="
$t.Cell(19,3).Range.Text = "50-3="
$t.Cell(19,4).Range.Text = "38+29="
$t.Cell(19,5).Range.Text = "90-41="
$t.Cell(20,1).Range.Text = "57-36="
$t.Cell(20,2).Range.Text = "96-75="
$t.Cell(20,3).Range.Text = "86-10="
$t.Cell(20,4).Range.Text = "21+18="
$t.Cell(20,5).Range.Text = "21+45="
